# Auto-generated script to apply 2025-05-07 crime data update
# Updates column L (year 2025 running total) and a couple of corrected
# 2023 (column J) values across the "Citywide Totals", "By Neighborhood"
# and individual neighborhood worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 2124
$ws.Range('L3').Value = 2135
$ws.Range('J4').Value = 1865
$ws.Range('L4').Value = 593
$ws.Range('L6').Value = 1919
$ws.Range('J7').Value = 29337
$ws.Range('L7').Value = 6894

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range('L2').Value = 7
$ws.Range('L7').Value = 17

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L2').Value = 123
$ws.Range('L3').Value = 149
$ws.Range('L6').Value = 113
$ws.Range('L7').Value = 433

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('L2').Value = 52
$ws.Range('L3').Value = 65
$ws.Range('L7').Value = 162

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L2').Value = 81
$ws.Range('L6').Value = 108
$ws.Range('L7').Value = 312

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L3').Value = 71
$ws.Range('L6').Value = 72
$ws.Range('L7').Value = 245

$ws = $wb.Worksheets.Item('New City')
$ws.Range('L3').Value = 42
$ws.Range('L7').Value = 132

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('L4').Value = 11
$ws.Range('L6').Value = 22
$ws.Range('L7').Value = 107

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L8').Value = 433
$ws.Range('L9').Value = 43
$ws.Range('L12').Value = 17
$ws.Range('L15').Value = 47
$ws.Range('L19').Value = 193
$ws.Range('L20').Value = 179
$ws.Range('L27').Value = 70
$ws.Range('J29').Value = 1553
$ws.Range('L29').Value = 357
$ws.Range('L33').Value = 312
$ws.Range('L34').Value = 43
$ws.Range('L36').Value = 100
$ws.Range('L37').Value = 245
$ws.Range('L41').Value = 33
$ws.Range('L42').Value = 215
$ws.Range('L43').Value = 54
$ws.Range('L44').Value = 49
$ws.Range('L45').Value = 12
$ws.Range('L48').Value = 96
$ws.Range('L63').Value = 21
$ws.Range('L64').Value = 50
$ws.Range('L65').Value = 132
$ws.Range('L66').Value = 15
$ws.Range('L67').Value = 250
$ws.Range('L69').Value = 17
$ws.Range('L72').Value = 30
$ws.Range('L76').Value = 72
$ws.Range('L78').Value = 95
$ws.Range('L79').Value = 188
$ws.Range('L83').Value = 162
$ws.Range('L85').Value = 363
$ws.Range('L89').Value = 90
$ws.Range('L91').Value = 100
$ws.Range('L93').Value = 38
$ws.Range('L94').Value = 80
$ws.Range('L97').Value = 63
$ws.Range('L99').Value = 107
$ws.Range('J101').Value = 29337
$ws.Range('L101').Value = 6894

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L2').Value = 73
$ws.Range('L7').Value = 250

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L3').Value = 126
$ws.Range('J4').Value = 85
$ws.Range('J7').Value = 1553
$ws.Range('L7').Value = 357

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('L2').Value = 13
$ws.Range('L4').Value = 27
$ws.Range('L6').Value = 35
$ws.Range('L7').Value = 96

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('L2').Value = 63
$ws.Range('L7').Value = 193

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('L2').Value = 19
$ws.Range('L7').Value = 49

$ws = $wb.Worksheets.Item('River North')
$ws.Range('L6').Value = 33
$ws.Range('L7').Value = 72

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('L6').Value = 8
$ws.Range('L7').Value = 33

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('L4').Value = 23
$ws.Range('L6').Value = 72
$ws.Range('L7').Value = 215

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('L6').Value = 28
$ws.Range('L7').Value = 95

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('L3').Value = 34
$ws.Range('L6').Value = 16
$ws.Range('L7').Value = 100

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('L3').Value = 65
$ws.Range('L6').Value = 39
$ws.Range('L7').Value = 188

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('L6').Value = 14
$ws.Range('L7').Value = 50

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('L4').Value = 15
$ws.Range('L7').Value = 179

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('L2').Value = 44
$ws.Range('L7').Value = 100

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('L3').Value = 10
$ws.Range('L7').Value = 38

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('L3').Value = 13
$ws.Range('L7').Value = 43

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('L6').Value = 24
$ws.Range('L7').Value = 80

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('L6').Value = 9
$ws.Range('L7').Value = 47

$ws = $wb.Worksheets.Item('North Center')
$ws.Range('L6').Value = 5
$ws.Range('L7').Value = 15

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('L3').Value = 19
$ws.Range('L6').Value = 11
$ws.Range('L7').Value = 43

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('L6').Value = 39
$ws.Range('L7').Value = 63

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('L6').Value = 22
$ws.Range('L7').Value = 90

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('L2').Value = 17
$ws.Range('L7').Value = 70

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('L6').Value = 18
$ws.Range('L7').Value = 54

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L2').Value = 112
$ws.Range('L7').Value = 363

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('L2').Value = 13
$ws.Range('L7').Value = 30

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range('L6').Value = 3
$ws.Range('L7').Value = 12

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('L6').Value = 4
$ws.Range('L7').Value = 17
